$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations")

$ws.Range("A7").Value = "Dog"
$ws.Range("B7").Value = "Dog"
$ws.Range("C7").Value = "Hund"
$ws.Range("D7").Value = "Need review"
